# Applies the "view section & progress increase added" edit:
# each data sheet gets a new header row (column names) inserted at row 1,
# pushing the existing data rows down by one; the users_have_courses sheet
# also gains five new blank (but styled) rows at the bottom; a handful of
# column widths are set, and the active sheet/selection moves from
# users_have_courses (sheet9) to courses (sheet1).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 9: users_have_courses  -> users_id, courses_id, status, rating, progress
# ---------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item(9)
$ws9.Rows("1:1").Insert()
$ws9.Range("A1").Value = "users_id"
$ws9.Range("B1").Value = "courses_id"
$ws9.Range("C1").Value = "status"
$ws9.Range("D1").Value = "rating"
$ws9.Range("E1").Value = "progress"

# Five new blank styled rows (10-14), matching format of the existing data rows.
$ws9.Range("A2:E4").Copy()
$ws9.Range("A10:E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws9.Range("E13").Select()

# ---------------------------------------------------------------------
# Sheet 8: users -> id, email, password, first_name, last_name, role,
#                   verification_token, is_verified, is_approved
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)
$ws8.Rows("1:1").Insert()
$ws8.Range("B1").Value = "id"
$ws8.Range("C1").Value = "email"
$ws8.Range("D1").Value = "password"
$ws8.Range("E1").Value = "first_name"
$ws8.Range("F1").Value = "last_name"
$ws8.Range("G1").Value = "role"
$ws8.Range("H1").Value = "verification_token"
$ws8.Range("I1").Value = "is_verified"
$ws8.Range("J1").Value = "is_approved"

$ws8.Columns.Item(5).ColumnWidth = 11.631510416666666
$ws8.Columns.Item(6).ColumnWidth = 16.365885416666668
$ws8.Columns.Item(8).ColumnWidth = 15.096354166666666
$ws8.Columns.Item(10).ColumnWidth = 11.233072916666666

$ws8.Activate()
$ws8.Range("B1:J6").Select()

# ---------------------------------------------------------------------
# Sheet 7: teachers -> users_id, phone_number, linked_in_account
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Rows("1:1").Insert()
$ws7.Range("A1").Value = "users_id"
$ws7.Range("B1").Value = "phone_number"
$ws7.Range("C1").Value = "linked_in_account"

$ws7.Columns.Item(3).ColumnWidth = 31.631510416666668

$ws7.Range("A1:C2").Select()

# ---------------------------------------------------------------------
# Sheet 6: tags -> id, expertise_area
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows("1:1").Insert()
$ws6.Range("A1").Value = "id"
$ws6.Range("B1").Value = "expertise_area"

$ws6.Columns.Item(2).ColumnWidth = 14.631510416666666

$ws6.Range("A1:B2").Select()

# ---------------------------------------------------------------------
# Sheet 5: sections -> id, title, content, description, external_link, courses_id
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows("1:1").Insert()
$ws5.Range("A1").Value = "id"
$ws5.Range("B1").Value = "title"
$ws5.Range("C1").Value = "content"
$ws5.Range("D1").Value = "description"
$ws5.Range("E1").Value = "external_link"
$ws5.Range("F1").Value = "courses_id"

$ws5.Columns.Item(4).ColumnWidth = 13.764322916666666
$ws5.Columns.Item(5).ColumnWidth = 13.166666666666666
$ws5.Columns.Item(6).ColumnWidth = 10.166666666666666

$ws5.Range("A1:F6").Select()

# ---------------------------------------------------------------------
# Sheet 4: objectives -> id, description
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows("1:1").Insert()
$ws4.Range("A1").Value = "id"
$ws4.Range("B1").Value = "description"

$ws4.Columns.Item(2).ColumnWidth = 16.166666666666668

$ws4.Range("A1:B3").Select()

# ---------------------------------------------------------------------
# Sheet 3: courses_have_tags -> courses_id, tags_id
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows("1:1").Insert()
$ws3.Range("A1").Value = "courses_id"
$ws3.Range("B1").Value = "tags_id"

$ws3.Range("A1:B4").Select()

# ---------------------------------------------------------------------
# Sheet 2: courses_have_objectives -> objectives_id, courses_id
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows("1:1").Insert()
$ws2.Range("A1").Value = "objectives_id"
$ws2.Range("B1").Value = "courses_id"

$ws2.Columns.Item(1).ColumnWidth = 13.096354166666666

$ws2.Range("A1:B4").Select()

# ---------------------------------------------------------------------
# Sheet 1: courses -> id, title, description, home_page_pic, owner_id,
#                      is_active, is_premium, course_rating
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("1:1").Insert()
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "title"
$ws1.Range("C1").Value = "description"
$ws1.Range("D1").Value = "home_page_pic"
$ws1.Range("E1").Value = "owner_id"
$ws1.Range("F1").Value = "is_active"
$ws1.Range("G1").Value = "is_premium"
$ws1.Range("H1").Value = "course_rating"

$ws1.Columns.Item(3).ColumnWidth = 11.365885416666666
$ws1.Columns.Item(4).ColumnWidth = 12.897135416666666
$ws1.Columns.Item(7).ColumnWidth = 9.432291666666666
$ws1.Columns.Item(8).ColumnWidth = 10.764322916666666

# Make courses the active sheet / tab, with H7 selected, as the final state.
$ws1.Activate()
$ws1.Range("H7").Select()
